# Development Log.xlsx edit
# Commit: printNumberInTileBag() - add a new dev-log entry row above the
# existing one, recording the latest progress; old entry shifts down a row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dev Log")

# --- 0. Remove the stale mailto hyperlink that lived on C4 (the old time
#        value). The replacement data is a plain number, not a contact link.
$ws.Range("C4").Hyperlinks.Delete()

# --- 1. Push the existing data row (row 4) down into row 5, carrying its
#        values AND formatting with it, before we overwrite row 4.
$ws.Range("B4:F4").Copy($ws.Range("B5:F5"))

# --- 2. Write the new log entry into row 4.
$ws.Range("B4").Value = 45359
$ws.Range("C4").Value = 20.32
$ws.Range("D4").Value = "Game now displays human, and or computer tiles (if open game)."
$ws.Range("E4").Value = "Focus on getting methods working in a  'rough and ready' status."
$ws.Range("F4").Value = "With 14 days remaining, I'm becoming increasing aware of how quickly time is running out. Focus now on brute-forcibly trying to make the game work, with refinement and unit-testing pushed to the lesser background for now.  Becoming more au fait with ArrayList and Map data structures and how best to implement them.`nDeveloping an agile/scrum mental mindset, proritizing a 'product log' of fixes, encompassing the 'bigger picture', with lots of short, tactical scrum-style sprint-log periods."

# --- 3. Both the new row-4 value and the carried-down row-5 value are plain
#        numbers now (hours.minutes notation), not a time-of-day serial, so
#        switch their number format from time to General, and re-key the
#        carried-down value the same way (13 h 26 m -> 13.26).
$ws.Range("C4:C5").NumberFormat = "General"
$ws.Range("C5").Value = 13.26

# --- 4. The grid around the header/data rows gains a top rule on each cell
#        (so row 3 separates from row 2, row 4 from row 3, and row 5 from
#        row 4) in addition to the pre-existing thin left/right/bottom.
$ws.Range("B3:F5").Borders.Item(8).LineStyle = 1
$ws.Range("B3:F5").Borders.Item(8).Weight = 2
$ws.Range("B3:F5").Borders.Item(8).ColorIndex = 1

# --- 5. Conditional formatting used to cover only B4:F4; now that there are
#        two data rows it needs to cover B5:F5 too, each driven by its own
#        row's column-A flag.
$existing = $ws.Range("B4").FormatConditions.Item(1)
$fc5 = $ws.Range("B5:F5").FormatConditions
$rule5 = $fc5.Add(2, 0, '=$A5=1')
$rule5.Font.Bold = $true
$rule5.Font.Italic = $false
$rule5.Font.Color = $existing.Font.Color

# --- 6. Row 2's thick bottom rule is gone now that row 3/4/5 form a taller
#        contiguous block; drop the extra weight so it matches its normal
#        medium box border.
$ws.Range("B2:F2").Borders.Item(9).Weight = -4138
